# Updated cryptos list on Thu Sep 26 11:31:08 UTC 2024 with GitHub Actions
#
# A new coin ("OKB") entered the ranking at position 41 (0-based), pushing
# every coin below it down by one row and dropping the last row
# (BabyDogeCoin) off the bottom of the 50-row table. All Price/Volume(1h)
# values were also refreshed with the latest snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into column D (Price). Most prices look like plain
# numbers ("28.08", "0.999", ...) but the sheet stores them as literal text
# (e.g. thousand-dot-separated "64.365.66", or subscript-notation
# "0.0\u20830912"), so any value that *would* be auto-coerced into a real
# number by Excel's Value setter must be forced back to text with a leading
# apostrophe (classic "store as text" trick) to preserve the exact digits/
# formatting. Values that already fail numeric parsing (extra dots,
# subscript glyphs, etc.) are left alone since Excel keeps them as text anyway.
function Set-DValue($sheet, $row, $val) {
    $ref = "D" + $row
    if ($val -match '^-?[0-9]+(\.[0-9]+)?$') {
        $sheet.Range($ref).Value = "'" + $val
    } else {
        $sheet.Range($ref).Value = $val
    }
}

# --- Insert the new "OKB" row at position 43, shifting Filecoin..BabyDogeCoin
# --- down by one, then drop what is now the trailing duplicate 52nd row
# --- (the old BabyDogeCoin row) so the table stays at 51 rows (header + 50).
$ws.Rows("43:43").Insert()
$ws.Rows("52:52").Delete()

# Row 43 is brand new (blank) after the insert - copy the numbering column's
# look (border/bold/center style) from the row above, then fill its contents.
$ws.Range("A42").Copy()
$ws.Range("A43").PasteSpecial(-4122)
$ws.Range("A43").Value = 41

$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"

# Column A is a plain positional counter (0-based row index), not tied to a
# particular coin, so `Rows.Insert` dragging 41..49 down into 44..52 needs
# undoing: re-pin rows 44..51 back to the same 42..49 sequence they always
# had (the diff leaves every A-cell untouched).
for ($i = 44; $i -le 51; $i++) {
    $ws.Range("A" + $i).Value = $i - 2
}

# --- Refresh Price (D) / Volume(1h) (E) for every data row 2..51 -----------
Set-DValue $ws 2 '64.365.66'
$ws.Range('E2').Value = '  +1.22%  '
Set-DValue $ws 3 '2.625.51'
$ws.Range('E3').Value = '  +0.22%  '
Set-DValue $ws 4 '1.00'
$ws.Range('E4').Value = '  +0.09%  '
Set-DValue $ws 5 '595.54'
$ws.Range('E5').Value = '  +0.14%  '
Set-DValue $ws 6 '152.65'
$ws.Range('E6').Value = '  +1.55%  '
Set-DValue $ws 7 '1.00'
$ws.Range('E7').Value = '  +0.07%  '
Set-DValue $ws 8 '0.590'
$ws.Range('E8').Value = '  +0.43%  '
Set-DValue $ws 9 '0.114'
$ws.Range('E9').Value = '  +4.88%  '
Set-DValue $ws 10 '5.84'
$ws.Range('E10').Value = '  +2.48%  '
Set-DValue $ws 11 '0.394'
$ws.Range('E11').Value = '  +3.39%  '
Set-DValue $ws 12 '0.152'
$ws.Range('E12').Value = '  +1.17%  '
Set-DValue $ws 13 '28.08'
$ws.Range('E13').Value = '  +1.50%  '
Set-DValue $ws 14 '3.100.96'
$ws.Range('E14').Value = '  +0.39%  '
Set-DValue $ws 15 '0.0000172'
$ws.Range('E15').Value = '  +13.74%  '
Set-DValue $ws 16 '64.338.92'
$ws.Range('E16').Value = '  +1.47%  '
Set-DValue $ws 17 '2.655.88'
$ws.Range('E17').Value = '  +1.12%  '
Set-DValue $ws 18 '12.30'
$ws.Range('E18').Value = '  -0.27%  '
Set-DValue $ws 19 '4.77'
$ws.Range('E19').Value = '  +2.55%  '
Set-DValue $ws 20 '349.73'
$ws.Range('E20').Value = '  +0.78%  '
Set-DValue $ws 21 '7.07'
$ws.Range('E21').Value = '  +2.99%  '
Set-DValue $ws 22 '1.00'
$ws.Range('E22').Value = '  +0.23%  '
Set-DValue $ws 23 '67.74'
$ws.Range('E23').Value = '  +2.20%  '
Set-DValue $ws 24 '1.70'
$ws.Range('E24').Value = '  -1.75%  '
Set-DValue $ws 25 '9.25'
$ws.Range('E25').Value = '  +0.29%  '
Set-DValue $ws 26 '1.67'
$ws.Range('E26').Value = '  -0.35%  '
Set-DValue $ws 27 '8.35'
$ws.Range('E27').Value = '  +1.65%  '
Set-DValue $ws 28 '550.10'
$ws.Range('E28').Value = '  -2.74%  '
Set-DValue $ws 29 '0.162'
$ws.Range('E29').Value = '  +0.44%  '
Set-DValue $ws 30 '0.999'
$ws.Range('E30').Value = '  -0.09%  '
Set-DValue $ws 31 '0.0₃0912'
$ws.Range('E31').Value = '  +7.85%  '
Set-DValue $ws 32 '2.08'
$ws.Range('E32').Value = '  +1.86%  '
Set-DValue $ws 33 '1.83'
$ws.Range('E33').Value = '  +4.47%  '
Set-DValue $ws 34 '5.52'
$ws.Range('E34').Value = '  +5.31%  '
Set-DValue $ws 35 '6.23'
$ws.Range('E35').Value = '  +2.07%  '
Set-DValue $ws 36 '0.421'
$ws.Range('E36').Value = '  +3.23%  '
Set-DValue $ws 37 '165.61'
$ws.Range('E37').Value = '  -1.72%  '
Set-DValue $ws 38 '20.10'
$ws.Range('E38').Value = '  +3.73%  '
Set-DValue $ws 39 '2.00'
$ws.Range('E39').Value = '  +3.19%  '
Set-DValue $ws 40 '1.00'
$ws.Range('E40').Value = '  +0.10%  '
Set-DValue $ws 41 '0.999'
$ws.Range('E41').Value = '  +0.01%  '
Set-DValue $ws 42 '168.65'
$ws.Range('E42').Value = '  +1.04%  '
Set-DValue $ws 43 '42.36'
$ws.Range('E43').Value = '  +6.17%  '
Set-DValue $ws 44 '4.11'
$ws.Range('E44').Value = '  +5.09%  '
Set-DValue $ws 45 '23.13'
$ws.Range('E45').Value = '  +7.74%  '
Set-DValue $ws 46 '2.24'
$ws.Range('E46').Value = '  +12.81%  '
Set-DValue $ws 47 '0.0590'
$ws.Range('E47').Value = '  -1.10%  '
Set-DValue $ws 48 '0.641'
$ws.Range('E48').Value = '  +2.06%  '
Set-DValue $ws 49 '0.0253'
$ws.Range('E49').Value = '  +1.42%  '
Set-DValue $ws 50 '0.0977'
$ws.Range('E50').Value = '  +1.56%  '
Set-DValue $ws 51 '19.35'
$ws.Range('E51').Value = '  +0.10%  '
